# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N ("Late") on the
# "Repayment schedule" sheet, shifting the old N/O/P columns to O/P/Q.
# Then make "Repayment schedule" the active sheet/tab (it was "Summary"
# before), and update its selection to J17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a blank column before column N (shifts N->O, O->P, P->Q)
$ws.Columns("N:N").Insert()

# Make the "Repayment schedule" sheet the active tab (was "Summary")
$ws.Activate()

# Update the selected cell on the now-active sheet
[void]$ws.Range("J17").Select()
